# "Error fixed in FD opening"
# Adds EquitableMortRegNo / EMRDate / Agentcode / DSACode / VehicleNo / RC.No
# columns (I:N) with sample data to the AccountOpening_Loan_ACOPL_TwoWh sheet,
# replacing the old ReferenceBy / Canvassername columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AccountOpening_Loan_ACOPL_TwoWh")
$ws.Activate()

# xlPasteFormats
$xlPasteFormats = -4122

# --- Column I: EquitableMortRegNo / Ad353 ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial($xlPasteFormats)
$ws.Range("I1").Value = "EquitableMortRegNo"

$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial($xlPasteFormats)
$ws.Range("I2").Value = "Ad353"

# --- Column J: EMRDate / 15/04/2025 ---
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial($xlPasteFormats)
$ws.Range("J1").Value = "EMRDate"

$ws.Range("H2").Copy()
$ws.Range("J2").PasteSpecial($xlPasteFormats)
$ws.Range("J2").Value = "15/04/2025"

# --- Column K: Agentcode / 7 ---
$ws.Range("H1").Copy()
$ws.Range("K1").PasteSpecial($xlPasteFormats)
$ws.Range("K1").Value = "Agentcode"

$ws.Range("G2").Copy()
$ws.Range("K2").PasteSpecial($xlPasteFormats)
$ws.Range("K2").Value = 7

# --- Column L: DSACode / 1 ---
$ws.Range("H1").Copy()
$ws.Range("L1").PasteSpecial($xlPasteFormats)
$ws.Range("L1").Value = "DSACode"

$ws.Range("G2").Copy()
$ws.Range("L2").PasteSpecial($xlPasteFormats)
$ws.Range("L2").Value = 1

# --- Column M: VehicleNo / TN 01 AB 1234 ---
$ws.Range("H1").Copy()
$ws.Range("M1").PasteSpecial($xlPasteFormats)
$ws.Range("M1").Value = "VehicleNo"

$ws.Range("G2").Copy()
$ws.Range("M2").PasteSpecial($xlPasteFormats)
$ws.Range("M2").Value = "TN 01 AB 1234"

# --- Column N: RC.No / 1234567890 ---
$ws.Range("H1").Copy()
$ws.Range("N1").PasteSpecial($xlPasteFormats)
$ws.Range("N1").Value = "RC.No"

$ws.Range("G2").Copy()
$ws.Range("N2").PasteSpecial($xlPasteFormats)
$ws.Range("N2").Value = 1234567890

$excel.CutCopyMode = $false

# Column widths for the newly-introduced columns (and J's revised width)
$ws.Columns.Item(9).ColumnWidth = 18.5
$ws.Columns.Item(10).ColumnWidth = 13.666666666666666
$ws.Columns.Item(11).ColumnWidth = 11.666666666666666
$ws.Columns.Item(13).ColumnWidth = 11.833333333333332
$ws.Columns.Item(14).ColumnWidth = 12.166666666666666

# Scroll the view over and select the new last header cell, like the author did
[void]$ws.Range("N1").Select()
